$p = $ppt.ActivePresentation

# --- 1) Update the "Update automatically" date placeholder text on the
#        Slide Master and on every Slide Layout: "1.01.2026" -> "2.01.2026".
function Update-DateShape($shape) {
    if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "1.01.2026") {
            $tr.Text = "2.01.2026"
        }
    }
}

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape $master.Shapes.Item($i)
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape $layout.Shapes.Item($i)
    }
}

# --- 2) Slide 2 ("Usuwanie elementow / dequeue"): Decision shape condition
#        "size > 1" -> "size >= 1" (single run edit, same run formatting).
$slide2 = $p.Slides.Item(2)
$dec2 = $slide2.Shapes.Item(6)
$tr2 = $dec2.TextFrame.TextRange
$sub2 = $tr2.Characters(5, $tr2.Length - 4)
$sub2.Text = " >= 1"

# --- 3) Slide 3 ("Sprawdzanie elementu z przodu kolejki / peek"): same
#        condition, but split into a new run for " >= " ahead of "1".
$slide3 = $p.Slides.Item(3)
$dec3 = $slide3.Shapes.Item(5)
$tr3 = $dec3.TextFrame.TextRange
$sub3 = $tr3.Characters(5, $tr3.Length - 5)
$sub3.Font.Size = 9
$sub3.Text = " >= "
